# edit c_click None issue
# Updates a few item quantities that were left wrong/empty ("None") by the
# c_click handler bug, across the three affected sheets.

$wb = $excel.ActiveWorkbook

# --- 매점판매 (sheet 2): 삼겹살 qty 4 -> 5, 김밥 qty 7 -> 9 ---
$wsSnack = $wb.Worksheets.Item("매점판매")
$wsSnack.Range("C4").Value = 5
$wsSnack.Range("C7").Value = 9

# --- 장의용품 (sheet 3): 맥주 qty 11 -> 10 ---
$wsFuneral = $wb.Worksheets.Item("장의용품")
$wsFuneral.Range("C5").Value = 10
# touch column D (benign, pre-existing style) so the sheet's used-range
# picks up the same click-through artifact seen in the recorded edit
$wsFuneral.Range("D1").Style = $wsFuneral.Range("A1").Style

# --- 기타 (sheet 5): 대패삼겹살 qty 43 -> 42, 치킨 qty 1 -> 0 ---
$wsEtc = $wb.Worksheets.Item("기타")
$wsEtc.Range("C9").Value = 42
$wsEtc.Range("C13").Value = 0
